$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.388023
$ws.Cells.Item(2, 8).Value = 7.164069
$ws.Cells.Item(2, 9).Value = 0.629429112239379
$ws.Cells.Item(2, 10).Value = 0.629429112239379
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.051093
$ws.Cells.Item(2, 14).Value = 0.153279
$ws.Cells.Item(2, 15).Value = 0.01450579975525089
$ws.Cells.Item(2, 16).Value = 0.01450579975525089
$ws.Cells.Item(2, 17).Value = 0.122011259139
$ws.Cells.Item(2, 18).Value = 1.098101332251
$ws.Cells.Item(2, 19).Value = 0.009130372662269766
$ws.Cells.Item(2, 20).Value = 0.00913037266226977
$ws.Cells.Item(3, 7).Value = 2.388023
$ws.Cells.Item(3, 8).Value = 7.164069
$ws.Cells.Item(3, 9).Value = 0.629429112239379
$ws.Cells.Item(3, 10).Value = 0.629429112239379
$ws.Cells.Item(3, 15).Value = 0.2313022967634575
$ws.Cells.Item(3, 16).Value = 0.2313022967634575
$ws.Cells.Item(3, 17).Value = 1.945531094184333
$ws.Cells.Item(3, 18).Value = 17.509779847659
$ws.Cells.Item(3, 19).Value = 0.1455883993107524
$ws.Cells.Item(3, 20).Value = 0.1455883993107525
$ws.Cells.Item(4, 7).Value = 2.388023
$ws.Cells.Item(4, 8).Value = 7.164069
$ws.Cells.Item(4, 9).Value = 0.629429112239379
$ws.Cells.Item(4, 10).Value = 0.629429112239379
$ws.Cells.Item(4, 13).Value = 2.656449666666667
$ws.Cells.Item(4, 14).Value = 7.969348999999999
$ws.Cells.Item(4, 15).Value = 0.7541919034812916
$ws.Cells.Item(4, 16).Value = 0.7541919034812917
$ws.Cells.Item(4, 17).Value = 6.343662902342333
$ws.Cells.Item(4, 18).Value = 57.092966121081
$ws.Cells.Item(4, 19).Value = 0.4747103402663568
$ws.Cells.Item(4, 20).Value = 0.4747103402663569
$ws.Cells.Item(5, 9).Value = 0.05397838594281493
$ws.Cells.Item(5, 10).Value = 0.05397838594281493
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.051093
$ws.Cells.Item(5, 14).Value = 0.153279
$ws.Cells.Item(5, 15).Value = 0.01450579975525089
$ws.Cells.Item(5, 16).Value = 0.01450579975525089
$ws.Cells.Item(5, 17).Value = 0.010463403594
$ws.Cells.Item(5, 18).Value = 0.094170632346
$ws.Cells.Item(5, 19).Value = 0.0007829996575981227
$ws.Cells.Item(5, 20).Value = 0.0007829996575981228
$ws.Cells.Item(6, 9).Value = 0.05397838594281493
$ws.Cells.Item(6, 10).Value = 0.05397838594281493
$ws.Cells.Item(6, 15).Value = 0.2313022967634575
$ws.Cells.Item(6, 16).Value = 0.2313022967634575
$ws.Cells.Item(6, 19).Value = 0.01248532464415742
$ws.Cells.Item(6, 20).Value = 0.01248532464415742
$ws.Cells.Item(7, 9).Value = 0.05397838594281493
$ws.Cells.Item(7, 10).Value = 0.05397838594281493
$ws.Cells.Item(7, 13).Value = 2.656449666666667
$ws.Cells.Item(7, 14).Value = 7.969348999999999
$ws.Cells.Item(7, 15).Value = 0.7541919034812916
$ws.Cells.Item(7, 16).Value = 0.7541919034812917
$ws.Cells.Item(7, 17).Value = 0.5440178691695555
$ws.Cells.Item(7, 18).Value = 4.896160822525999
$ws.Cells.Item(7, 19).Value = 0.04071006164105939
$ws.Cells.Item(7, 20).Value = 0.04071006164105939
$ws.Cells.Item(8, 7).Value = 1.201136333333333
$ws.Cells.Item(8, 8).Value = 3.603409
$ws.Cells.Item(8, 9).Value = 0.3165925018178061
$ws.Cells.Item(8, 10).Value = 0.3165925018178061
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.051093
$ws.Cells.Item(8, 14).Value = 0.153279
$ws.Cells.Item(8, 15).Value = 0.01450579975525089
$ws.Cells.Item(8, 16).Value = 0.01450579975525089
$ws.Cells.Item(8, 17).Value = 0.061369658679
$ws.Cells.Item(8, 18).Value = 0.552326928111
$ws.Cells.Item(8, 19).Value = 0.004592427435382998
$ws.Cells.Item(8, 20).Value = 0.004592427435382999
$ws.Cells.Item(9, 7).Value = 1.201136333333333
$ws.Cells.Item(9, 8).Value = 3.603409
$ws.Cells.Item(9, 9).Value = 0.3165925018178061
$ws.Cells.Item(9, 10).Value = 0.3165925018178061
$ws.Cells.Item(9, 15).Value = 0.2313022967634575
$ws.Cells.Item(9, 16).Value = 0.2313022967634575
$ws.Cells.Item(9, 17).Value = 0.9785701749332222
$ws.Cells.Item(9, 18).Value = 8.807131574399
$ws.Cells.Item(9, 19).Value = 0.07322857280854764
$ws.Cells.Item(9, 20).Value = 0.07322857280854765
$ws.Cells.Item(10, 7).Value = 1.201136333333333
$ws.Cells.Item(10, 8).Value = 3.603409
$ws.Cells.Item(10, 9).Value = 0.3165925018178061
$ws.Cells.Item(10, 10).Value = 0.3165925018178061
$ws.Cells.Item(10, 13).Value = 2.656449666666667
$ws.Cells.Item(10, 14).Value = 7.969348999999999
$ws.Cells.Item(10, 15).Value = 0.7541919034812916
$ws.Cells.Item(10, 16).Value = 0.7541919034812917
$ws.Cells.Item(10, 17).Value = 3.190758212304555
$ws.Cells.Item(10, 18).Value = 28.716823910741
$ws.Cells.Item(10, 19).Value = 0.2387715015738755
$ws.Cells.Item(10, 20).Value = 0.2387715015738755